$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column from 2023-11-13 (45243)
# to 2023-11-14 (45244) for rows 2 through 7.
foreach ($row in 2..7) {
    $ws.Range("C$row").Value = 45244
}
